$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = 10.25574820157856
$ws.Range("D2").Value = 8.975944711274391
$ws.Range("E2").Value = 13.70060417621827
$ws.Range("F2").Value = 34.4933259223417
$ws.Range("G2").Value = 38.12757742536402
$ws.Range("H2").Value = 16.14600372607444
$ws.Range("I2").Value = 26.73130812648908
$ws.Range("J2").Value = 10.16132451069928
$ws.Range("L2").Value = 10.32893753701926
$ws.Range("O2").Value = 25.95845761623582

$ws.Range("C3").Value = 10.25419415395164
$ws.Range("D3").Value = 8.98551772318287
$ws.Range("E3").Value = 13.71301761476361
$ws.Range("F3").Value = 34.32215378162802
$ws.Range("G3").Value = 37.69371076815759
$ws.Range("H3").Value = 16.13245442039494
$ws.Range("I3").Value = 26.63979269713989
$ws.Range("J3").Value = 10.17838943921144
$ws.Range("L3").Value = 10.34618474666638
$ws.Range("O3").Value = 25.85418634298545

$ws.Range("C4").Value = 10.25496736024086
$ws.Range("D4").Value = 8.992546182110372
$ws.Range("E4").Value = 13.72265027229662
$ws.Range("F4").Value = 34.22483470244256
$ws.Range("G4").Value = 37.43552600402915
$ws.Range("H4").Value = 16.12714361120178
$ws.Range("I4").Value = 26.58940240184983
$ws.Range("J4").Value = 10.18997233868108
$ws.Range("L4").Value = 10.35787407113777
$ws.Range("O4").Value = 25.79586358366984

$ws.Range("C5").Value = 10.25571774317546
$ws.Range("D5").Value = 8.995699989109648
$ws.Range("E5").Value = 13.72708169799857
$ws.Range("F5").Value = 34.18716484512647
$ws.Range("G5").Value = 37.33250096491635
$ws.Range("H5").Value = 16.12573820316198
$ws.Range("I5").Value = 26.57034141363606
$ws.Range("J5").Value = 10.19497065615731
$ws.Range("L5").Value = 10.36291439552206
$ws.Range("O5").Value = 25.77354764553505

$ws.Range("C6").Value = 10.2558686579161
$ws.Range("D6").Value = 8.996241179479803
$ws.Range("E6").Value = 13.72784810662954
$ws.Range("F6").Value = 34.18103076345412
$ws.Range("G6").Value = 37.31552933438708
$ws.Range("H6").Value = 16.12555072288737
$ws.Range("I6").Value = 26.56726574579979
$ws.Range("J6").Value = 10.19581743484617
$ws.Range("L6").Value = 10.36376806775679
$ws.Range("O6").Value = 25.76993022197648

$ws.Range("C7").Value = 10.25497571658052
$ws.Range("D7").Value = 8.992587542252114
$ws.Range("E7").Value = 13.72270798663793
$ws.Range("F7").Value = 34.2243185826257
$ws.Range("G7").Value = 37.43412755747855
$ws.Range("H7").Value = 16.12712158243222
$ws.Range("I7").Value = 26.58913935436005
$ws.Range("J7").Value = 10.1900386209079
$ws.Range("L7").Value = 10.35794092537174
$ws.Range("O7").Value = 25.79555672589973

$ws.Range("C8").Value = 10.25485445262201
$ws.Range("D8").Value = 8.979006887707557
$ws.Range("E8").Value = 13.7044672763149
$ws.Range("F8").Value = 34.43271069899574
$ws.Range("G8").Value = 37.97635093970705
$ws.Range("H8").Value = 16.14070899972143
$ws.Range("I8").Value = 26.69855913196371
$ws.Range("J8").Value = 10.16697936907392
$ws.Range("L8").Value = 10.33465639849387
$ws.Range("O8").Value = 25.92133296682629

$ws.Range("C9").Value = 10.26827549925821
$ws.Range("D9").Value = 8.961489428108825
$ws.Range("E9").Value = 13.68463057890936
$ws.Range("F9").Value = 34.90157944826323
$ws.Range("G9").Value = 39.09920595849613
$ws.Range("H9").Value = 16.19110787293874
$ws.Range("I9").Value = 26.95842681131263
$ws.Range("J9").Value = 10.13051262704167
$ws.Range("L9").Value = 10.29770340712711
$ws.Range("O9").Value = 26.21234824143723

$ws.Range("C10").Value = 10.28638961953528
$ws.Range("D10").Value = 8.954153750359788
$ws.Range("E10").Value = 13.67973727800888
$ws.Range("F10").Value = 35.28053749958851
$ws.Range("G10").Value = 39.95276507461053
$ws.Range("H10").Value = 16.24243436187499
$ws.Range("I10").Value = 27.17588296160218
$ws.Range("J10").Value = 10.10903585112754
$ws.Range("L10").Value = 10.27584109295462
$ws.Range("O10").Value = 26.45193824389055

$ws.Range("C11").Value = 10.29640291300809
$ws.Range("D11").Value = 8.952012810803087
$ws.Range("E11").Value = 13.67960411940153
$ws.Range("F11").Value = 35.45989777936907
$ws.Range("G11").Value = 40.34560918664414
$ws.Range("H11").Value = 16.26884059538297
$ws.Range("I11").Value = 27.28030770749336
$ws.Range("J11").Value = 10.10041518649371
$ws.Range("L11").Value = 10.26703866062764
$ws.Range("O11").Value = 26.56622826055489

$ws.Range("C12").Value = 10.30044764780768
$ws.Range("D12").Value = 8.951373534015241
$ws.Range("E12").Value = 13.67985363855853
$ws.Range("F12").Value = 35.52877180817477
$ws.Range("G12").Value = 40.49488415481667
$ws.Range("H12").Value = 16.2792747404914
$ws.Range("I12").Value = 27.32061752924849
$ws.Range("J12").Value = 10.09731563409763
$ws.Range("L12").Value = 10.26386934274443
$ws.Range("O12").Value = 26.61024105061272

$ws.Range("C13").Value = 10.29956532796505
$ws.Range("D13").Value = 8.951503597632751
$ws.Range("E13").Value = 13.67978657735517
$ws.Range("F13").Value = 35.51389690145573
$ws.Range("G13").Value = 40.46271452596221
$ws.Range("H13").Value = 16.27700831258685
$ws.Range("I13").Value = 27.31190241010126
$ws.Range("J13").Value = 10.0979758496479
$ws.Range("L13").Value = 10.26454462532401
$ws.Range("O13").Value = 26.60072994580705

$ws.Range("C14").Value = 10.29673061937078
$ws.Range("D14").Value = 8.951956784298828
$ws.Range("E14").Value = 13.67961864235948
$ws.Range("F14").Value = 35.46554524596029
$ws.Range("G14").Value = 40.3578806450056
$ws.Range("H14").Value = 16.26969033741627
$ws.Range("I14").Value = 27.28360880729905
$ws.Range("J14").Value = 10.10015688154467
$ws.Range("L14").Value = 10.26677463483671
$ws.Range("O14").Value = 26.56983468329441

$ws.Range("C15").Value = 10.29502715164686
$ws.Range("D15").Value = 8.952256685074699
$ws.Range("E15").Value = 13.67955480695902
$ws.Range("F15").Value = 35.43605126447486
$ws.Range("G15").Value = 40.29372954679285
$ws.Range("H15").Value = 16.26526431818752
$ws.Range("I15").Value = 27.26637721108171
$ws.Range("J15").Value = 10.10151429156754
$ws.Range("L15").Value = 10.26816192297258
$ws.Range("O15").Value = 26.55100511959468

$ws.Range("C16").Value = 10.28577073563206
$ws.Range("D16").Value = 8.954317688712704
$ws.Range("E16").Value = 13.67978801038091
$ws.Range("F16").Value = 35.26895207073283
$ws.Range("G16").Value = 39.92717132315089
$ws.Range("H16").Value = 16.24076976162575
$ws.Range("I16").Value = 27.16916719476811
$ws.Range("J16").Value = 10.10962232749499
$ws.Range("L16").Value = 10.2764393184123
$ws.Range("O16").Value = 26.44457336149633

$ws.Range("C17").Value = 10.2805449594243
$ws.Range("D17").Value = 8.955888023426349
$ws.Range("E17").Value = 13.68046642217078
$ws.Range("F17").Value = 35.1681937341711
$ws.Range("G17").Value = 39.70336541921074
$ws.Range("H17").Value = 16.2265228486533
$ws.Range("I17").Value = 27.1109241412186
$ws.Range("J17").Value = 10.11489044382609
$ws.Range("L17").Value = 10.28180968445773
$ws.Range("O17").Value = 26.38061843306225

$ws.Range("C18").Value = 10.27770625809721
$ws.Range("D18").Value = 8.956903856558299
$ws.Range("E18").Value = 13.68105368136502
$ws.Range("F18").Value = 35.11089929227043
$ws.Range("G18").Value = 39.57507543001888
$ws.Range("H18").Value = 16.21861652075249
$ws.Range("I18").Value = 27.07794380719676
$ws.Range("J18").Value = 10.11802871762254
$ws.Range("L18").Value = 10.28500616954735
$ws.Range("O18").Value = 26.34433383519401

$ws.Range("C19").Value = 10.27677386884405
$ws.Range("D19").Value = 8.957267158325902
$ws.Range("E19").Value = 13.68128639011146
$ws.Range("F19").Value = 35.09161501483199
$ws.Range("G19").Value = 39.53171807743799
$ws.Range("H19").Value = 16.21598920358608
$ws.Range("I19").Value = 27.06686719854773
$ws.Range("J19").Value = 10.11910987589872
$ws.Range("L19").Value = 10.28610693461117
$ws.Range("O19").Value = 26.33213531316683

$ws.Range("C20").Value = 10.28108397978846
$ws.Range("D20").Value = 8.955709206869741
$ws.Range("E20").Value = 13.68037381653578
$ws.Range("F20").Value = 35.17885176091222
$ws.Range("G20").Value = 39.72714571647992
$ws.Range("H20").Value = 16.22800967227069
$ws.Range("I20").Value = 27.11707062021454
$ws.Range("J20").Value = 10.11431844882351
$ws.Range("L20").Value = 10.28122686701973
$ws.Range("O20").Value = 26.3873749370964

$ws.Range("C21").Value = 10.29755639520382
$ws.Range("D21").Value = 8.951819023872952
$ws.Range("E21").Value = 13.67965983716073
$ws.Range("F21").Value = 35.47972179541751
$ws.Range("G21").Value = 40.38866008521553
$ws.Range("H21").Value = 16.27182804914001
$ws.Range("I21").Value = 27.29189872801601
$ws.Range("J21").Value = 10.09951178696417
$ws.Range("L21").Value = 10.26611517981371
$ws.Range("O21").Value = 26.57888969274287

$ws.Range("C22").Value = 10.30979517933339
$ws.Range("D22").Value = 8.950275620908531
$ws.Range("E22").Value = 13.68094095816672
$ws.Range("F22").Value = 35.68189594445141
$ws.Range("G22").Value = 40.82392603959587
$ws.Range("H22").Value = 16.30299689102792
$ws.Range("I22").Value = 27.41061395752961
$ws.Range("J22").Value = 10.09079580500553
$ws.Range("L22").Value = 10.25719439537708
$ws.Range("O22").Value = 26.70831781823669

$ws.Range("C23").Value = 10.30312904882639
$ws.Range("D23").Value = 8.951008143761168
$ws.Range("E23").Value = 13.68009764378383
$ws.Range("F23").Value = 35.57350140083865
$ws.Range("G23").Value = 40.59139490750044
$ws.Range("H23").Value = 16.28613167351546
$ws.Range("I23").Value = 27.34685432140864
$ws.Range("J23").Value = 10.09535987393026
$ws.Range("L23").Value = 10.26186827498663
$ws.Range("O23").Value = 26.63885931430912

$ws.Range("C24").Value = 10.28083977250414
$ws.Range("D24").Value = 8.955789697735735
$ws.Range("E24").Value = 13.68041506911912
$ws.Range("F24").Value = 35.17403129195975
$ws.Range("G24").Value = 39.71639345259933
$ws.Range("H24").Value = 16.22733659296019
$ws.Range("I24").Value = 27.11429022403656
$ws.Range("J24").Value = 10.11457670645479
$ws.Range("L24").Value = 10.28149001926072
$ws.Range("O24").Value = 26.38431881195839

$ws.Range("C25").Value = 10.26318924705696
$ws.Range("D25").Value = 8.965254705995902
$ws.Range("E25").Value = 13.68829407869671
$ws.Range("F25").Value = 34.76851179398151
$ws.Range("G25").Value = 38.78985591801308
$ws.Range("H25").Value = 16.1749476266201
$ws.Range("I25").Value = 26.88338189682072
$ws.Range("J25").Value = 10.13944290318607
$ws.Range("L25").Value = 10.30677011373447
$ws.Range("O25").Value = 26.21234824143723
